# Glossary export update: re-shuffle the "Instrument Data" / "Auxiliary Data" /
# "Verification" rows, and strip the bold+boxed header styling back to the
# default (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-point the three affected term rows at their new definitions ---------
# row 3 (was "Instrument Data") -> becomes "Verification"
$ws.Range("A3").Value = "Verification"
$ws.Range("D3").Value = "Verification serves as a means to evaluate the reliability of the Data in the absence of a Reference dataset, allowing for an assessment of its standalone performance. It involves confirming the consistency and internal coherence of the Data without direct comparison to external Reference sources."
$ws.Range("G3").Value = "- KCEO"

# row 4 (was "Auxiliary Data") -> becomes "Instrument Data"
$ws.Range("A4").Value = "Instrument Data"
$ws.Range("D4").Value = "Data produced and transmitted by the science and engineering sensors of an instrument, and, in the spacecraft environment, any additional data packaged with the instrument’s sensor data by virtue of services provided"
$ws.Range("G4").Value = "- [EO Data Stewardship Glossary](https://ceos.org/document_management/Working_Groups/WGISS/Interest_Groups/Data_Stewardship/White_Papers/EO-DataStewardshipGlossary.pdf)"

# row 6 (was "Verification") -> becomes "Auxiliary Data"
$ws.Range("A6").Value = "Auxiliary Data"
$ws.Range("D6").Value = "The Data required for instrument processing, which does not originate in the instrument itself or from the satellite. Some auxiliary Data will be generated in the ground segment, whilst other Data will be provided from external sources."
$ws.Range("G6").Value = "- CEOS-ARD PFS template 20220302"

# --- Strip the header row's bold font + thin box border back to default ----
$ws.Range("A1:G1").ClearFormats()
